$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: num_customers 70 -> 71 (cohort_size D34 stays 2256, retention_rate recalculated)
$ws.Range("C34").Value = 71
$ws.Range("E34").Value = 0.03147163120567376

# Row 36: num_customers 112 -> 113 (cohort_size D36 stays 1930, retention_rate recalculated)
$ws.Range("C36").Value = 113
$ws.Range("E36").Value = 0.05854922279792746

# Row 37: num_customers 702 -> 707, cohort_size 702 -> 707 (retention_rate stays 1)
$ws.Range("C37").Value = 707
$ws.Range("D37").Value = 707
$ws.Range("E37").Value = 1
